$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 78
$ws.Range("F7").Value = 148
$ws.Range("F12").Value = 5574
$ws.Range("F13").Value = 75
$ws.Range("F14").Value = 6380
$ws.Range("F16").Value = 434
$ws.Range("F19").Value = 600
$ws.Range("F24").Value = 10297
$ws.Range("F25").Value = 1959
$ws.Range("F26").Value = 2184
$ws.Range("F29").Value = 2205
$ws.Range("F31").Value = 88
$ws.Range("F34").Value = 74
$ws.Range("F35").Value = 2131
$ws.Range("F39").Value = 5314

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 0
$ws.Range("F8").Value = 8
$ws.Range("F18").Value = 911
$ws.Range("F20").Value = 0
$ws.Range("F22").Value = 0

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 23

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 40
$ws.Range("F5").Value = 0
$ws.Range("F8").Value = 9152
$ws.Range("F9").Value = 148
$ws.Range("F13").Value = 323
$ws.Range("F18").Value = 6380
$ws.Range("F19").Value = 6380
$ws.Range("F21").Value = 434
$ws.Range("F23").Value = 600
$ws.Range("F25").Value = 213
$ws.Range("F26").Value = 161
$ws.Range("F28").Value = 10297
$ws.Range("F29").Value = 1959
$ws.Range("F30").Value = 2184
$ws.Range("F32").Value = 2205
$ws.Range("F33").Value = 84
$ws.Range("F37").Value = 2131
$ws.Range("F38").Value = 0
$ws.Range("F40").Value = 5314
$ws.Range("F41").Value = 1208
$ws.Range("F48").Value = 1392
$ws.Range("F50").Value = 0
